$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update revised F/G values for rows 334-365
$ws.Range("F334").Value = 204109 ; $ws.Range("G334").Value = 3483
$ws.Range("F335").Value = 130516 ; $ws.Range("G335").Value = 2982
$ws.Range("F336").Value = 101348 ; $ws.Range("G336").Value = 3300
$ws.Range("F337").Value = 103594 ; $ws.Range("G337").Value = 2890
$ws.Range("F338").Value = 224163 ; $ws.Range("G338").Value = 3120
$ws.Range("F339").Value = 653844 ; $ws.Range("G339").Value = 5585
$ws.Range("F340").Value = 379600 ; $ws.Range("G340").Value = 3263
$ws.Range("F341").Value = 294578 ; $ws.Range("G341").Value = 3675
$ws.Range("F342").Value = 179365 ; $ws.Range("G342").Value = 3082
$ws.Range("F343").Value = 133895 ; $ws.Range("G343").Value = 2960
$ws.Range("F344").Value = 136459
$ws.Range("F345").Value = 290423 ; $ws.Range("G345").Value = 3306
$ws.Range("F346").Value = 660704 ; $ws.Range("G346").Value = 4719
$ws.Range("F347").Value = 338911 ; $ws.Range("G347").Value = 2878
$ws.Range("F348").Value = 231263 ; $ws.Range("G348").Value = 3229
$ws.Range("F349").Value = 159655 ; $ws.Range("G349").Value = 2755
$ws.Range("F350").Value = 128050 ; $ws.Range("G350").Value = 2970
$ws.Range("F351").Value = 149955 ; $ws.Range("G351").Value = 2818
$ws.Range("F352").Value = 305941 ; $ws.Range("G352").Value = 3542
$ws.Range("F353").Value = 711389 ; $ws.Range("G353").Value = 5214
$ws.Range("F354").Value = 304175 ; $ws.Range("G354").Value = 2781
$ws.Range("F355").Value = 221475 ; $ws.Range("G355").Value = 3431
$ws.Range("F356").Value = 160254 ; $ws.Range("G356").Value = 2892
$ws.Range("F357").Value = 137689 ; $ws.Range("G357").Value = 2872
$ws.Range("F358").Value = 157825 ; $ws.Range("G358").Value = 2603
$ws.Range("F359").Value = 319570 ; $ws.Range("G359").Value = 3345
$ws.Range("F360").Value = 733895 ; $ws.Range("G360").Value = 5015
$ws.Range("F361").Value = 327377 ; $ws.Range("G361").Value = 2567
$ws.Range("F362").Value = 223662 ; $ws.Range("G362").Value = 3072
$ws.Range("F363").Value = 184834 ; $ws.Range("G363").Value = 2711
$ws.Range("F364").Value = 161731 ; $ws.Range("G364").Value = 2377
$ws.Range("F365").Value = 183115 ; $ws.Range("G365").Value = 2393

# Add new rows 366-368
$ws.Range("A366").Value = 44260
$ws.Range("B366").Value = 322104
$ws.Range("C366").Value = 15721
$ws.Range("D366").Value = 2522
$ws.Range("E366").Value = 7739
$ws.Range("F366").Value = 308946
$ws.Range("G366").Value = 2535

$ws.Range("A367").Value = 44261
$ws.Range("B367").Value = 323390
$ws.Range("C367").Value = 7611
$ws.Range("D367").Value = 1286
$ws.Range("E367").Value = 7836
$ws.Range("F367").Value = 639534
$ws.Range("G367").Value = 3263

$ws.Range("A368").Value = 44262
$ws.Range("B368").Value = 323786
$ws.Range("C368").Value = 2747
$ws.Range("D368").Value = 396
$ws.Range("E368").Value = 7921
$ws.Range("F368").Value = 267722
$ws.Range("G368").Value = 1789